$wb = $excel.ActiveWorkbook

# --- "input" sheet: replace the placeholder/test example row (row 2) with a
#     real-world Nepal example based on OSM online, and clear the no-longer
#     relevant shapefile columns. ("klein foutje bij analyses_indirect.py eruit!")
$wsInput = $wb.Worksheets.Item("input")

$wsInput.Range("A2").Value = "Nepal"
$wsInput.Range("D2").Value = "Network based on OSM online"
$wsInput.Range("E2").Value = "npl_admbnda_adm0_nd_20201117.shp"
$wsInput.Range("G2").ClearContents()
$wsInput.Range("H2").ClearContents()
$wsInput.Range("L2").Value = "drive"
$wsInput.Range("M2").Value = "motorway, trunk, primary, secondary"

# --- view/selection state: "input" becomes the active sheet/tab, with the
#     whole of row 2 selected (the example row that was just edited).
$wsInput.Activate()
$wsInput.Range("A2:XFD2").Select()

# --- "explanation" sheet loses the tabSelected flag and its selection moves
#     back to A2 (single cell, no longer the active tab).
$wsExplanation = $wb.Worksheets.Item("explanation")
$wsExplanation.Range("A2").Select()

# Re-activate "input" last so it remains the workbook's active sheet/tab.
$wsInput.Activate()
